$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.215.73"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "3.500.29"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.35"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.19"
$ws.Range("E6").Value = "  -2.23%  "

$ws.Range("D7").Value = "3.498.29"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.135"
$ws.Range("E10").Value = "  -4.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "8.00"
$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  -2.53%  "

$ws.Range("D13").Value = "4.093.89"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("E14").Value = "  -3.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.27"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").Value = "3.502.35"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "66.236.33"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.53"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  -3.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.65"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("E23").Value = "  -1.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.92"
$ws.Range("E24").Value = "  -2.32%  "

$ws.Range("D25").Value = "3.634.70"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.25"
$ws.Range("E28").Value = "  -5.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").Value = "  -3.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.46"
$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("E33").Value = "  -8.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.16"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").Value = "3.488.48"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.72"
$ws.Range("E37").Value = "  -3.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.60"
$ws.Range("E38").Value = "  -4.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.76"
$ws.Range("E39").Value = "  -2.93%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.08"
$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0856"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("E43").Value = "  -4.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.882"
$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.90"
$ws.Range("E45").Value = "  -9.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.32"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("E47").Value = "  -8.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.75"
$ws.Range("E48").Value = "  -11.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("E50").Value = "  -4.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.939"
$ws.Range("E51").Value = "  -3.47%  "
